$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -38
$ws.Range("H28").Value = 3610.158
$ws.Range("I28").Value = 4576.5356
$ws.Range("J28").Value = 904.3
$ws.Range("K28").Value = 4576.5356
$ws.Range("L28").Value = 904.3
$ws.Range("M28").Value = -4091.5356
$ws.Range("N28").Value = -1874.3
$ws.Range("H64").Value = 90912216
$ws.Range("I64").Value = 125002050
$ws.Range("J64").Value = 5986.6665
$ws.Range("K64").Value = 125002050
$ws.Range("L64").Value = 5986.6665
$ws.Range("M64").Value = -125001802
$ws.Range("N64").Value = -6482.6665
$ws.Range("H67").Value = 90912216
$ws.Range("I67").Value = 125002050
$ws.Range("J67").Value = 5986.6665
$ws.Range("K67").Value = 125002050
$ws.Range("L67").Value = 5986.6665
$ws.Range("M67").Value = -125001192
$ws.Range("N67").Value = -7702.6665
$ws.Range("H74").Value = 2944.1052
$ws.Range("J74").Value = 2996.5557
$ws.Range("L74").Value = 2996.5557
$ws.Range("N74").Value = -4868.5557
$ws.Range("H77").Value = 2944.1052
$ws.Range("J77").Value = 2996.5557
$ws.Range("L77").Value = 14982.7785
$ws.Range("N77").Value = -24342.7785
$ws.Range("H97").Value = 788.8889
$ws.Range("J97").Value = 788.8889
$ws.Range("L97").Value = 2366.6667
$ws.Range("N97").Value = -3358.6667
$ws.Range("H98").Value = 2086
$ws.Range("J98").Value = 5995
$ws.Range("L98").Value = 5995
$ws.Range("N98").Value = -8991
$ws.Range("H121").Value = 1087.1428
$ws.Range("J121").Value = 1242
$ws.Range("L121").Value = 3726
$ws.Range("N121").Value = -7220
$ws.Range("H122").Value = 2086
$ws.Range("J122").Value = 5995
$ws.Range("L122").Value = 17985
$ws.Range("N122").Value = -22885

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7429.25
$ws.Range("I74").Value = 1089.5
$ws.Range("J74").Value = 26448.5
$ws.Range("K74").Value = 1089.5
$ws.Range("L74").Value = 26448.5
$ws.Range("M74").Value = -215.5
$ws.Range("N74").Value = -28196.5
$ws.Range("H76").Value = 36923
$ws.Range("J76").Value = 36923
$ws.Range("L76").Value = 36923
$ws.Range("N76").Value = -37599
$ws.Range("H77").Value = 7429.25
$ws.Range("I77").Value = 1089.5
$ws.Range("J77").Value = 26448.5
$ws.Range("K77").Value = 5447.5
$ws.Range("L77").Value = 132242.5
$ws.Range("M77").Value = -1079.5
$ws.Range("N77").Value = -140978.5
$ws.Range("H79").Value = 36923
$ws.Range("J79").Value = 36923
$ws.Range("L79").Value = 36923
$ws.Range("N79").Value = -39263

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 5388.933
$ws.Range("I25").Value = 722.3333
$ws.Range("J25").Value = 8500
$ws.Range("K25").Value = 722.3333
$ws.Range("L25").Value = 8500
$ws.Range("M25").Value = -487.3333
$ws.Range("N25").Value = -8970
$ws.Range("H134").Value = 1253.0834
$ws.Range("I134").Value = 929
$ws.Range("J134").Value = 2225.3333
$ws.Range("K134").Value = 2787
$ws.Range("L134").Value = 6675.999899999999
$ws.Range("M134").Value = -252
$ws.Range("N134").Value = -11745.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H99").Value = 1702.1052
$ws.Range("I99").Value = 1458.3636
$ws.Range("K99").Value = 1458.3636
$ws.Range("M99").Value = 39.63640000000009
$ws.Range("H126").Value = 1702.1052
$ws.Range("I126").Value = 1458.3636
$ws.Range("K126").Value = 4375.0908
$ws.Range("M126").Value = -1905.0908

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 109
$ws.Range("I6").Value = 54.444443
$ws.Range("K6").Value = 163.333329
$ws.Range("M6").Value = -50.33332899999999
$ws.Range("H132").Value = 50723.8
$ws.Range("I132").Value = 667.7692
$ws.Range("J132").Value = 143685
$ws.Range("K132").Value = 6009.922799999999
$ws.Range("L132").Value = 1293165
$ws.Range("M132").Value = -3479.922799999999
$ws.Range("N132").Value = -1298225

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 62859.855
$ws.Range("I11").Value = 45001.5
$ws.Range("J11").Value = 70003.2
$ws.Range("K11").Value = 45001.5
$ws.Range("L11").Value = 70003.2
$ws.Range("M11").Value = -44862.5
$ws.Range("N11").Value = -70281.2
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H113").Value = 201544.8
$ws.Range("I113").Value = 334603.66
$ws.Range("J113").Value = 1956.5
$ws.Range("K113").Value = 334603.66
$ws.Range("L113").Value = 1956.5
$ws.Range("M113").Value = -332433.66
$ws.Range("N113").Value = -6296.5
$ws.Range("H122").Value = 1043.4445
$ws.Range("I122").Value = 984.5714
$ws.Range("J122").Value = 1249.5
$ws.Range("K122").Value = 2953.7142
$ws.Range("L122").Value = 3748.5
$ws.Range("M122").Value = -503.7142000000003
$ws.Range("N122").Value = -8648.5
$ws.Range("H126").Value = 1432.25
$ws.Range("I126").Value = 1160.2
$ws.Range("J126").Value = 1885.6666
$ws.Range("K126").Value = 3480.6
$ws.Range("L126").Value = 5656.9998
$ws.Range("M126").Value = -1010.6
$ws.Range("N126").Value = -10596.9998
$ws.Range("H132").Value = 45712.81
$ws.Range("I132").Value = 87664
$ws.Range("J132").Value = 3761.6155
$ws.Range("K132").Value = 262992
$ws.Range("L132").Value = 11284.8465
$ws.Range("M132").Value = -260462
$ws.Range("N132").Value = -16344.8465

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 44415.668
$ws.Range("J62").Value = 44415.668
$ws.Range("L62").Value = 44415.668
$ws.Range("N62").Value = -45663.668
$ws.Range("H65").Value = 44415.668
$ws.Range("J65").Value = 44415.668
$ws.Range("L65").Value = 133247.004
$ws.Range("N65").Value = -139487.004
$ws.Range("H68").Value = 1729.1
$ws.Range("I68").Value = 1769.3334
$ws.Range("J68").Value = 1668.75
$ws.Range("K68").Value = 1769.3334
$ws.Range("L68").Value = 1668.75
$ws.Range("M68").Value = -1020.3334
$ws.Range("N68").Value = -3166.75
$ws.Range("H71").Value = 1729.1
$ws.Range("I71").Value = 1769.3334
$ws.Range("J71").Value = 1668.75
$ws.Range("K71").Value = 8846.666999999999
$ws.Range("L71").Value = 8343.75
$ws.Range("M71").Value = -5102.666999999999
$ws.Range("N71").Value = -15831.75
$ws.Range("H122").Value = 41579.19
$ws.Range("I122").Value = 73985.71000000001
$ws.Range("J122").Value = 3771.5833
$ws.Range("K122").Value = 221957.13
$ws.Range("L122").Value = 11314.7499
$ws.Range("M122").Value = -219507.13
$ws.Range("N122").Value = -16214.7499

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 10095
$ws.Range("I70").Value = 10095
$ws.Range("K70").Value = 10095
$ws.Range("M70").Value = -9780
$ws.Range("H73").Value = 10095
$ws.Range("I73").Value = 10095
$ws.Range("K73").Value = 10095
$ws.Range("M73").Value = -9003
